$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)    # "About" sheet
$ws2 = $wb.Worksheets.Item(2)    # "PPEIdtIL" sheet

# ---------------------------------------------------------------------------
# "About" sheet: replace the old "Note:" paragraph (rows 10-14) with a new
# explanatory paragraph about the variable (rows 10-16), then re-add the
# original ACEEE-labeling paragraph further down (rows 18-22), leaving row 17
# blank.
# ---------------------------------------------------------------------------

$newNote = @(
  "This variable reflects improvement in efficiency components selected by consumers due",
  "to improved labeling. The labeling influences consumers who are buying appliances of all",
  "Quality levels, so it's represented as a simple percentage increase in the efficiency of",
  "components sold (at all quality levels). If Quality Levels are defined based on",
  "particular efficiency thresholds, this may mean the number of square feet served by",
  "components of a given quality level will not be accurate. It's just a question of the meaning",
  "of the labels given to each Quality Level."
)

$row = 10
foreach ($line in $newNote) {
  $ws1.Cells.Item($row, 1).Value = $line
  $row++
}

$oldAceeeParagraph = @(
  "The ACEEE study focused on ""appliance"" labeling.  We use the same percentage for labeling",
  "of heating equipment, as well as cooling and ventilation equipment, because they are",
  "similar (e.g. machines one buys in a store, which could readily bear labels, with similar",
  "costs and lifetimes as other major appliances).  We similarly assume the same rate",
  "applies to commercial and residential buildings."
)

$row = 18
foreach ($line in $oldAceeeParagraph) {
  $ws1.Cells.Item($row, 1).Value = $line
  $row++
}

$ws1.Range("A17:XFD30").Select()

# ---------------------------------------------------------------------------
# "PPEIdtIL" sheet: rename the header from "Building Component" to
# "Efficiency Improvement by Building Component (dimensionless)" and wrap the
# header row text, increasing its height.
# ---------------------------------------------------------------------------

$ws2.Range("A1").Value = "Efficiency Improvement by Building Component (dimensionless)"
$ws2.Range("A1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 45
